$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(3, "agriculture"),
    @(8, "tank"),
    @(16, "land_without_scrub"),
    @(21, "river"),
    @(30, "forest"),
    @(43, "land_without_scrub"),
    @(59, "road_n_railway"),
    @(68, "river"),
    @(69, "crop"),
    @(75, "road_n_railway"),
    @(81, "road_n_railway"),
    @(94, "river"),
    @(95, "agriculture"),
    @(102, "grassland"),
    @(123, "road_n_railway"),
    @(130, "road_n_railway"),
    @(141, "land_without_scrub"),
    @(144, "road_n_railway"),
    @(154, "crop"),
    @(158, "settlement"),
    @(162, "crop"),
    @(172, "grassland"),
    @(183, "crop"),
    @(185, "river"),
    @(188, "crop"),
    @(191, "road_n_railway"),
    @(198, "land_without_scrub"),
    @(204, "agriculture"),
    @(210, "agriculture"),
    @(218, "tank"),
    @(244, "river"),
    @(260, "crop"),
    @(276, "land_without_scrub"),
    @(282, "river"),
    @(286, "agriculture"),
    @(292, "river"),
    @(305, "agriculture"),
    @(309, "road_n_railway"),
    @(311, "settlement"),
    @(315, "grassland"),
    @(322, "grassland"),
    @(328, "river"),
    @(333, "river"),
    @(343, "land_without_scrub"),
    @(366, "road_n_railway"),
    @(373, "river"),
    @(392, "land_without_scrub"),
    @(394, "crop"),
    @(412, "crop"),
    @(416, "crop"),
    @(421, "agriculture"),
    @(429, "land_without_scrub"),
    @(432, "agriculture"),
    @(435, "forest"),
    @(441, "river"),
    @(456, "crop"),
    @(491, "land_without_scrub"),
    @(493, "crop"),
    @(498, "river"),
    @(502, "grassland"),
    @(516, "forest"),
    @(519, "forest"),
    @(526, "forest"),
    @(531, "river"),
    @(534, "agriculture"),
    @(542, "road_n_railway"),
    @(549, "river"),
    @(562, "crop"),
    @(568, "agriculture"),
    @(569, "tank"),
    @(594, "crop"),
    @(595, "road_n_railway"),
    @(598, "tank"),
    @(618, "grassland"),
    @(619, "land_without_scrub"),
    @(628, "grassland"),
    @(632, "land_without_scrub"),
    @(633, "agriculture"),
    @(634, "crop"),
    @(641, "road_n_railway"),
    @(646, "crop"),
    @(649, "land_without_scrub"),
    @(654, "river"),
    @(668, "river"),
    @(669, "road_n_railway"),
    @(671, "crop"),
    @(672, "land_without_scrub"),
    @(694, "road_n_railway"),
    @(695, "river"),
    @(696, "crop"),
    @(717, "river"),
    @(726, "land_without_scrub"),
    @(735, "grassland"),
    @(740, "river"),
    @(755, "river"),
    @(756, "grassland"),
    @(761, "land_without_scrub"),
    @(765, "land_without_scrub"),
    @(767, "river"),
    @(771, "crop"),
    @(780, "land_without_scrub"),
    @(788, "forest"),
    @(812, "grassland"),
    @(816, "river"),
    @(826, "river"),
    @(836, "crop"),
    @(843, "road_n_railway"),
    @(852, "river"),
    @(857, "river"),
    @(876, "land_without_scrub"),
    @(894, "land_without_scrub"),
    @(895, "land_without_scrub"),
    @(901, "forest"),
    @(911, "grassland"),
    @(913, "forest"),
    @(927, "agriculture"),
    @(934, "crop"),
    @(940, "crop"),
    @(945, "settlement"),
    @(948, "grassland"),
    @(952, "river"),
    @(962, "land_without_scrub"),
    @(965, "river"),
    @(967, "road_n_railway"),
    @(971, "river"),
    @(976, "grassland"),
    @(999, "land_without_scrub"),
    @(1003, "river"),
    @(1020, "settlement"),
    @(1021, "river"),
    @(1042, "road_n_railway"),
    @(1048, "crop"),
    @(1051, "river"),
    @(1069, "land_without_scrub"),
    @(1077, "river"),
    @(1091, "river"),
    @(1094, "crop"),
    @(1097, "settlement"),
    @(1099, "grassland"),
    @(1102, "tank"),
    @(1110, "land_without_scrub"),
    @(1129, "settlement"),
    @(1131, "grassland"),
    @(1137, "crop"),
    @(1152, "land_without_scrub"),
    @(1158, "land_without_scrub"),
    @(1169, "river"),
    @(1170, "river"),
    @(1172, "grassland"),
    @(1189, "river"),
    @(1190, "river"),
    @(1200, "river"),
    @(1201, "land_without_scrub"),
    @(1209, "grassland"),
    @(1213, "river"),
    @(1250, "river"),
    @(1251, "river"),
    @(1291, "river"),
    @(1294, "road_n_railway"),
    @(1297, "land_without_scrub"),
    @(1309, "river"),
    @(1320, "road_n_railway"),
    @(1321, "land_without_scrub"),
    @(1325, "grassland"),
    @(1326, "tank"),
    @(1338, "road_n_railway"),
    @(1356, "tank"),
    @(1357, "road_n_railway"),
    @(1363, "agriculture"),
    @(1369, "crop"),
    @(1374, "land_without_scrub"),
    @(1379, "land_without_scrub"),
    @(1385, "land_without_scrub"),
    @(1392, "grassland"),
    @(1419, "land_without_scrub"),
    @(1420, "road_n_railway"),
    @(1428, "crop"),
    @(1430, "road_n_railway"),
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], 2).Value = $u[1]
}